$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: B26, C26 get values "Yes" / "Marcelo" with Calibri 12 font (style 7)
$ws.Range("B26").Value = "Yes"
$ws.Range("B26").Font.Name = "Calibri"
$ws.Range("B26").Font.Size = 12

$ws.Range("C26").Value = "Marcelo"
$ws.Range("C26").Font.Name = "Calibri"
$ws.Range("C26").Font.Size = 12

# Row 38: B38, C38 get values "Yes" / "Javier" with Calibri 12 font (style 7)
$ws.Range("B38").Value = "Yes"
$ws.Range("B38").Font.Name = "Calibri"
$ws.Range("B38").Font.Size = 12

$ws.Range("C38").Value = "Javier"
$ws.Range("C38").Font.Name = "Calibri"
$ws.Range("C38").Font.Size = 12

# Row 43: B43 gets value "Yes" with Calibri 12 font (style 7) (C43 already has Javier)
$ws.Range("B43").Value = "Yes"
$ws.Range("B43").Font.Name = "Calibri"
$ws.Range("B43").Font.Size = 12

# Row 45: A45 style change only (keep value) to Arial 12 (style 9); C45 gets value "Bianca"
$ws.Range("A45").Font.Name = "Arial"
$ws.Range("A45").Font.Size = 12

$ws.Range("C45").Value = "Bianca"

# Row 46: C46 changes from Javier to Marcelo
$ws.Range("C46").Value = "Marcelo"
